# Update "想去人数" (F column) values on the "展览" sheet and the
# aggregated "全部类型" sheet, per the upstream gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 6-21) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 38
$wsExhibit.Range("F7").Value  = 593
$wsExhibit.Range("F8").Value  = 114
$wsExhibit.Range("F9").Value  = 8755
$wsExhibit.Range("F10").Value = 811
$wsExhibit.Range("F11").Value = 331
$wsExhibit.Range("F12").Value = 1145
$wsExhibit.Range("F13").Value = 989
$wsExhibit.Range("F14").Value = 112
$wsExhibit.Range("F15").Value = 47
$wsExhibit.Range("F17").Value = 236
$wsExhibit.Range("F18").Value = 262
$wsExhibit.Range("F19").Value = 67
$wsExhibit.Range("F20").Value = 231
$wsExhibit.Range("F21").Value = 1036

# --- Sheet "全部类型" (rows 7-23, offset by +1 vs "展览") ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 38
$wsAll.Range("F9").Value  = 593
$wsAll.Range("F10").Value = 114
$wsAll.Range("F11").Value = 8755
$wsAll.Range("F12").Value = 811
$wsAll.Range("F13").Value = 331
$wsAll.Range("F14").Value = 1145
$wsAll.Range("F15").Value = 989
$wsAll.Range("F16").Value = 112
$wsAll.Range("F17").Value = 47
$wsAll.Range("F19").Value = 236
$wsAll.Range("F20").Value = 262
$wsAll.Range("F21").Value = 67
$wsAll.Range("F22").Value = 231
$wsAll.Range("F23").Value = 1036
